$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.022.19"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.300.04"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'300.05"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'97.83"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").Value = "'0.519"
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'36.07"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "2.659.55"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "2.299.82"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "'0.787"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "42.910.71"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'12.76"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'68.23"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "'237.55"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'24.94"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -13.34%  "
$ws.Range("D30").Value = "'9.14"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'163.48"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").Value = "'33.00"
$ws.Range("E32").Value = "  -4.35%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("D35").Value = "'18.04"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").Value = "'4.77"
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'0.0696"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "'2.79"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").Value = "2.019.59"
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.26"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0286"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "'10.37"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'17.43"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'2.83"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "'54.29"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "2.525.43"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  -1.39%  "
